# Generate Report for Handoff
# Re-running the handoff report generation refreshed the timestamps for the
# 1e053c64-... / 35b44c5d-... entries and set the "ht" (handoff type)
# priority on the rows that were (re)handed off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows whose handoff priority flips from blank to "ht" in both locale sheets.
$rows = @(7, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}

# Refresh "Latest Handoff Datetime" on the zh-cn sheet for the same rows.
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-08-17 20:19:59"
}

# Refresh "Latest HO Xliff Generate Date" on the Overview sheet for the
# matching rows.
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-17 20:20:17"
}
